# Add three new worksheets (RoiProfiles, RoiCorners, RoiCenters) at the end of
# the workbook, each holding the same header row: shapes, label, description, image.

$wb = $excel.ActiveWorkbook

$sheetNames = @("RoiProfiles", "RoiCorners", "RoiCenters")
$headers = @("shapes", "label", "description", "image")

foreach ($name in $sheetNames) {
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $ws = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
    $ws.Name = $name

    for ($i = 0; $i -lt $headers.Length; $i++) {
        $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
    }
}
